$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("780÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "106÷4=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("541÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "266÷6=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("694÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "940÷7=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("166÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "647÷9=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("437÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "871÷3=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("327÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "830÷2=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("203÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "304÷9=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("792÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "837÷4=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("266÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "460÷5=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("369÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "958÷8=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("568÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "664÷2=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("261÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "146÷5=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("125÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "901÷9=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("734÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "681÷6=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("435÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "545÷2=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("463÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "569÷2=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("512÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "745÷5=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("985÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "809÷5=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("296÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "148÷9=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("648÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "280÷9=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("951÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "414÷9=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("547÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "901÷6=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("904÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "317÷3=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("961÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷6=", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("737÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "579÷7=", 2) | Out-Null
